$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 117
$ws1.Range("G2").Value = 50
$ws1.Range("F4").Value = 979

# Sheet "全部类型" (All types) - sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 117
$ws4.Range("G2").Value = 50
$ws4.Range("F4").Value = 979
